$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 791
$ws1.Range("F5").Value = 877
$ws1.Range("F6").Value = 2152
$ws1.Range("F7").Value = 191

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 791
$ws4.Range("F7").Value = 877
$ws4.Range("F8").Value = 2152
$ws4.Range("F10").Value = 191
